# Applies the 2024-05-08 violent crime data update across all affected sheets.
# Generated from the canonical OOXML diff: updates 2024 (column K) totals, and a
# handful of small 2022/2023 (columns I/J) corrections, for 46 neighborhoods plus
# the Citywide Totals and By Neighborhood rollup sheets.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 2560
$ws.Range('I3').Value = 7487
$ws.Range('K3').Value = 2473
$ws.Range('I4').Value = 1787
$ws.Range('K4').Value = 515
$ws.Range('K5').Value = 163
$ws.Range('K6').Value = 3078
$ws.Range('I7').Value = 26240
$ws.Range('K7').Value = 8789

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 64
$ws.Range('K7').Value = 261
$ws.Range('J8').Value = 1852
$ws.Range('K8').Value = 585
$ws.Range('K9').Value = 33
$ws.Range('K10').Value = 49
$ws.Range('K11').Value = 186
$ws.Range('K14').Value = 52
$ws.Range('K15').Value = 88
$ws.Range('K16').Value = 27
$ws.Range('K18').Value = 59
$ws.Range('K19').Value = 257
$ws.Range('K20').Value = 199
$ws.Range('K23').Value = 79
$ws.Range('K27').Value = 94
$ws.Range('K29').Value = 448
$ws.Range('K31').Value = 102
$ws.Range('K33').Value = 349
$ws.Range('K37').Value = 283
$ws.Range('K44').Value = 85
$ws.Range('K46').Value = 19
$ws.Range('K49').Value = 58
$ws.Range('K50').Value = 55
$ws.Range('K51').Value = 96
$ws.Range('K52').Value = 239
$ws.Range('K53').Value = 130
$ws.Range('K54').Value = 163
$ws.Range('K55').Value = 95
$ws.Range('K57').Value = 24
$ws.Range('J63').Value = 100
$ws.Range('K63').Value = 35
$ws.Range('K64').Value = 57
$ws.Range('K65').Value = 206
$ws.Range('I67').Value = 980
$ws.Range('K67').Value = 346
$ws.Range('K68').Value = 23
$ws.Range('K73').Value = 86
$ws.Range('K74').Value = 12
$ws.Range('K75').Value = 34
$ws.Range('K76').Value = 128
$ws.Range('K77').Value = 61
$ws.Range('K78').Value = 122
$ws.Range('K79').Value = 227
$ws.Range('K83').Value = 195
$ws.Range('K85').Value = 422
$ws.Range('K89').Value = 115
$ws.Range('K95').Value = 140
$ws.Range('K96').Value = 120
$ws.Range('I99').Value = 447
$ws.Range('K99').Value = 157
$ws.Range('I101').Value = 26240
$ws.Range('K101').Value = 8789

# Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K2').Value = 22
$ws.Range('K3').Value = 9
$ws.Range('K4').Value = 3
$ws.Range('K7').Value = 52

# West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 42
$ws.Range('K6').Value = 55
$ws.Range('K7').Value = 120

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K3').Value = 83
$ws.Range('K6').Value = 72
$ws.Range('K7').Value = 261

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 57
$ws.Range('K6').Value = 74
$ws.Range('K7').Value = 186

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 28
$ws.Range('K7').Value = 115

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 156
$ws.Range('K3').Value = 144
$ws.Range('K6').Value = 95
$ws.Range('K7').Value = 422

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 56
$ws.Range('K6').Value = 100
$ws.Range('K7').Value = 239

# Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K6').Value = 69
$ws.Range('K7').Value = 130

# Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 171
$ws.Range('K3').Value = 175
$ws.Range('J4').Value = 96
$ws.Range('K4').Value = 32
$ws.Range('K6').Value = 194
$ws.Range('J7').Value = 1852
$ws.Range('K7').Value = 585

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 74
$ws.Range('K3').Value = 59
$ws.Range('K6').Value = 45
$ws.Range('K7').Value = 195

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 98
$ws.Range('K3').Value = 127
$ws.Range('K6').Value = 96
$ws.Range('K7').Value = 349

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K3').Value = 46
$ws.Range('K7').Value = 140

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 74
$ws.Range('K3').Value = 99
$ws.Range('K6').Value = 86
$ws.Range('K7').Value = 283

# New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K6').Value = 84
$ws.Range('K7').Value = 206

# Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I3').Value = 158
$ws.Range('K3').Value = 58
$ws.Range('K4').Value = 8
$ws.Range('I7').Value = 447
$ws.Range('K7').Value = 157

# Gage Park
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 33
$ws.Range('K6').Value = 41
$ws.Range('K7').Value = 102

# North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 109
$ws.Range('K3').Value = 108
$ws.Range('I4').Value = 54
$ws.Range('K6').Value = 104
$ws.Range('I7').Value = 980
$ws.Range('K7').Value = 346

# Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K3').Value = 10
$ws.Range('K7').Value = 58

# Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K2').Value = 32
$ws.Range('K7').Value = 163

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 122
$ws.Range('K3').Value = 149
$ws.Range('K4').Value = 26
$ws.Range('K6').Value = 141
$ws.Range('K7').Value = 448

# Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K3').Value = 19
$ws.Range('K6').Value = 56

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 82
$ws.Range('K7').Value = 257

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K3').Value = 24
$ws.Range('K7').Value = 85

# River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K2').Value = 23
$ws.Range('K3').Value = 21
$ws.Range('K7').Value = 128

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K2').Value = 16
$ws.Range('K7').Value = 49

# Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 36
$ws.Range('K3').Value = 29
$ws.Range('K4').Value = 9
$ws.Range('K7').Value = 122

# Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K6').Value = 33
$ws.Range('K7').Value = 95

# Jefferson Park
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('K6').Value = 8
$ws.Range('K7').Value = 19

# Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K2').Value = 27
$ws.Range('K7').Value = 79

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 82
$ws.Range('K7').Value = 227

# Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K2').Value = 13
$ws.Range('K7').Value = 57

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K6').Value = 69
$ws.Range('K7').Value = 199

# Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K6').Value = 13
$ws.Range('K7').Value = 59

# Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K2').Value = 31
$ws.Range('K7').Value = 88

# Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K3').Value = 6
$ws.Range('K7').Value = 55

# Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K3').Value = 12
$ws.Range('K7').Value = 33

# Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 86

# Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 19
$ws.Range('K3').Value = 18
$ws.Range('K7').Value = 64

# Edgewater
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K2').Value = 25
$ws.Range('K7').Value = 94

# Pullman
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('K2').Value = 14
$ws.Range('K7').Value = 34

# Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 96

# North Park
$ws = $wb.Worksheets.Item('North Park')
$ws.Range('K6').Value = 7
$ws.Range('K7').Value = 23

# Mckinley Park
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K3').Value = 4
$ws.Range('K7').Value = 24

# Riverdale
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 30
$ws.Range('K7').Value = 61

# Bucktown
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('K6').Value = 17
$ws.Range('K7').Value = 27

# Printers Row
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('K2').Value = 1
$ws.Range('K7').Value = 12
